$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 70
$ws.Cells.Item(70, 8).Value = 2018.1818
$ws.Cells.Item(70, 9).Value = 1439.8
$ws.Cells.Item(70, 10).Value = 2500.1667
$ws.Cells.Item(70, 11).Value = 4319.4
$ws.Cells.Item(70, 12).Value = 7500.500100000001
$ws.Cells.Item(70, 13).Value = -4049.4
$ws.Cells.Item(70, 14).Value = -8040.500100000001
# ALC row 73
$ws.Cells.Item(73, 8).Value = 2018.1818
$ws.Cells.Item(73, 9).Value = 1439.8
$ws.Cells.Item(73, 10).Value = 2500.1667
$ws.Cells.Item(73, 11).Value = 4319.4
$ws.Cells.Item(73, 12).Value = 7500.500100000001
$ws.Cells.Item(73, 13).Value = -3383.4
$ws.Cells.Item(73, 14).Value = -9372.500100000001
# ALC row 93
$ws.Cells.Item(93, 8).Value = 23483.5
$ws.Cells.Item(93, 10).Value = 23483.5
$ws.Cells.Item(93, 12).Value = 23483.5
$ws.Cells.Item(93, 14).Value = -28475.5
# ALC row 98
$ws.Cells.Item(98, 8).Value = 331816.28
$ws.Cells.Item(98, 9).Value = 363656.6
$ws.Cells.Item(98, 10).Value = 2800
$ws.Cells.Item(98, 11).Value = 363656.6
$ws.Cells.Item(98, 12).Value = 2800
$ws.Cells.Item(98, 13).Value = -362158.6
$ws.Cells.Item(98, 14).Value = -5796
# ALC row 121
$ws.Cells.Item(121, 8).Value = 899.2857
$ws.Cells.Item(121, 10).Value = 1119
$ws.Cells.Item(121, 12).Value = 3357
$ws.Cells.Item(121, 14).Value = -6851
# ALC row 122
$ws.Cells.Item(122, 8).Value = 331816.28
$ws.Cells.Item(122, 9).Value = 363656.6
$ws.Cells.Item(122, 10).Value = 2800
$ws.Cells.Item(122, 11).Value = 1090969.8
$ws.Cells.Item(122, 12).Value = 8400
$ws.Cells.Item(122, 13).Value = -1088519.8
$ws.Cells.Item(122, 14).Value = -13300
# ALC row 132
$ws.Cells.Item(132, 8).Value = 27572.436
$ws.Cells.Item(132, 9).Value = 32064.031
$ws.Cells.Item(132, 10).Value = 2868.6667
$ws.Cells.Item(132, 11).Value = 96192.09299999999
$ws.Cells.Item(132, 12).Value = 8606.000100000001
$ws.Cells.Item(132, 13).Value = -93662.09299999999
$ws.Cells.Item(132, 14).Value = -13666.0001

$ws = $wb.Worksheets.Item("ARM")
# ARM row 10
$ws.Cells.Item(10, 8).Value = 300002000
$ws.Cells.Item(10, 9).Value = 900000000
$ws.Cells.Item(10, 11).Value = 900000000
$ws.Cells.Item(10, 13).Value = -899999830
# ARM row 32
$ws.Cells.Item(32, 8).Value = 25992.229
$ws.Cells.Item(32, 9).Value = 5296.114
$ws.Cells.Item(32, 11).Value = 5296.114
$ws.Cells.Item(32, 13).Value = -5009.114
# ARM row 35
$ws.Cells.Item(35, 8).Value = 1737
$ws.Cells.Item(35, 9).Value = 1737
$ws.Cells.Item(35, 11).Value = 1737
$ws.Cells.Item(35, 13).Value = -1331
# ARM row 63
$ws.Cells.Item(63, 8).Value = 8856.286
$ws.Cells.Item(63, 10).Value = 9267.5
$ws.Cells.Item(63, 12).Value = 9267.5
$ws.Cells.Item(63, 14).Value = -10639.5
# ARM row 66
$ws.Cells.Item(66, 8).Value = 8856.286
$ws.Cells.Item(66, 10).Value = 9267.5
$ws.Cells.Item(66, 12).Value = 46337.5
$ws.Cells.Item(66, 14).Value = -53201.5
# ARM row 74
$ws.Cells.Item(74, 8).Value = 3392.0962
$ws.Cells.Item(74, 9).Value = 927.37836
$ws.Cells.Item(74, 10).Value = 9471.733
$ws.Cells.Item(74, 11).Value = 927.37836
$ws.Cells.Item(74, 12).Value = 9471.733
$ws.Cells.Item(74, 13).Value = -53.37836000000004
$ws.Cells.Item(74, 14).Value = -11219.733
# ARM row 77
$ws.Cells.Item(77, 8).Value = 3392.0962
$ws.Cells.Item(77, 9).Value = 927.37836
$ws.Cells.Item(77, 10).Value = 9471.733
$ws.Cells.Item(77, 11).Value = 4636.8918
$ws.Cells.Item(77, 12).Value = 47358.665
$ws.Cells.Item(77, 13).Value = -268.8918000000003
$ws.Cells.Item(77, 14).Value = -56094.665
# ARM row 132
$ws.Cells.Item(132, 8).Value = 2433.5193
$ws.Cells.Item(132, 9).Value = 1954.3954
$ws.Cells.Item(132, 10).Value = 4722.6665
$ws.Cells.Item(132, 11).Value = 5863.1862
$ws.Cells.Item(132, 12).Value = 14167.9995
$ws.Cells.Item(132, 13).Value = -3333.1862
$ws.Cells.Item(132, 14).Value = -19227.9995

$ws = $wb.Worksheets.Item("BSM")
# BSM row 22
$ws.Cells.Item(22, 8).Value = 246.16667
$ws.Cells.Item(22, 9).Value = 215.6
$ws.Cells.Item(22, 10).Value = 399
$ws.Cells.Item(22, 11).Value = 215.6
$ws.Cells.Item(22, 12).Value = 399
$ws.Cells.Item(22, 13).Value = -42.59999999999999
$ws.Cells.Item(22, 14).Value = -745
# BSM row 24
$ws.Cells.Item(24, 8).Value = 11254.2
$ws.Cells.Item(24, 9).Value = 12817.75
$ws.Cells.Item(24, 10).Value = 5000
$ws.Cells.Item(24, 11).Value = 12817.75
$ws.Cells.Item(24, 12).Value = 5000
$ws.Cells.Item(24, 13).Value = -12582.75
$ws.Cells.Item(24, 14).Value = -5470
# BSM row 99
$ws.Cells.Item(99, 8).Value = 1931
$ws.Cells.Item(99, 9).Value = 1489.6
$ws.Cells.Item(99, 11).Value = 1489.6
$ws.Cells.Item(99, 13).Value = 8.400000000000091
# BSM row 138
$ws.Cells.Item(138, 8).Value = 50000
$ws.Cells.Item(138, 10).Value = 50000
$ws.Cells.Item(138, 12).Value = 50000
$ws.Cells.Item(138, 14).Value = -60280

$ws = $wb.Worksheets.Item("CRP")
# CRP row 31
$ws.Cells.Item(31, 8).Value = 5326.2
$ws.Cells.Item(31, 9).Value = 1677.5667
$ws.Cells.Item(31, 10).Value = 10799.15
$ws.Cells.Item(31, 11).Value = 1677.5667
$ws.Cells.Item(31, 12).Value = 10799.15
$ws.Cells.Item(31, 13).Value = -1382.5667
$ws.Cells.Item(31, 14).Value = -11389.15
# CRP row 34
$ws.Cells.Item(34, 8).Value = 5326.2
$ws.Cells.Item(34, 9).Value = 1677.5667
$ws.Cells.Item(34, 10).Value = 10799.15
$ws.Cells.Item(34, 11).Value = 1677.5667
$ws.Cells.Item(34, 12).Value = 10799.15
$ws.Cells.Item(34, 13).Value = -1475.5667
$ws.Cells.Item(34, 14).Value = -11203.15
# CRP row 93
$ws.Cells.Item(93, 8).Value = 14690.3
$ws.Cells.Item(93, 9).Value = 12487.875
$ws.Cells.Item(93, 10).Value = 23500
$ws.Cells.Item(93, 11).Value = 12487.875
$ws.Cells.Item(93, 12).Value = 23500
$ws.Cells.Item(93, 13).Value = -10615.875
$ws.Cells.Item(93, 14).Value = -27244
# CRP row 99
$ws.Cells.Item(99, 8).Value = 2340.8928
$ws.Cells.Item(99, 9).Value = 1817.421
$ws.Cells.Item(99, 10).Value = 3446
$ws.Cells.Item(99, 11).Value = 1817.421
$ws.Cells.Item(99, 12).Value = 3446
$ws.Cells.Item(99, 13).Value = -319.421
$ws.Cells.Item(99, 14).Value = -6442
# CRP row 103
$ws.Cells.Item(103, 8).Value = 10219
$ws.Cells.Item(103, 9).Value = 6821.7144
$ws.Cells.Item(103, 10).Value = 34000
$ws.Cells.Item(103, 11).Value = 6821.7144
$ws.Cells.Item(103, 12).Value = 34000
$ws.Cells.Item(103, 13).Value = -5649.7144
$ws.Cells.Item(103, 14).Value = -36344
# CRP row 126
$ws.Cells.Item(126, 8).Value = 2340.8928
$ws.Cells.Item(126, 9).Value = 1817.421
$ws.Cells.Item(126, 10).Value = 3446
$ws.Cells.Item(126, 11).Value = 5452.263
$ws.Cells.Item(126, 12).Value = 10338
$ws.Cells.Item(126, 13).Value = -2982.263
$ws.Cells.Item(126, 14).Value = -15278
# CRP row 132
$ws.Cells.Item(132, 8).Value = 1652.579
$ws.Cells.Item(132, 9).Value = 1174.1212
$ws.Cells.Item(132, 10).Value = 4810.4
$ws.Cells.Item(132, 11).Value = 3522.3636
$ws.Cells.Item(132, 12).Value = 14431.2
$ws.Cells.Item(132, 13).Value = -992.3636000000001
$ws.Cells.Item(132, 14).Value = -19491.2

$ws = $wb.Worksheets.Item("CUL")
# CUL row 122
$ws.Cells.Item(122, 8).Value = 1262.2273
$ws.Cells.Item(122, 9).Value = 386.875
$ws.Cells.Item(122, 10).Value = 1762.4286
$ws.Cells.Item(122, 11).Value = 3481.875
$ws.Cells.Item(122, 12).Value = 15861.8574
$ws.Cells.Item(122, 13).Value = -1031.875
$ws.Cells.Item(122, 14).Value = -20761.8574

$ws = $wb.Worksheets.Item("GSM")
# GSM row 122
$ws.Cells.Item(122, 8).Value = 2039.8
$ws.Cells.Item(122, 9).Value = 900
$ws.Cells.Item(122, 10).Value = 2324.75
$ws.Cells.Item(122, 11).Value = 2700
$ws.Cells.Item(122, 12).Value = 6974.25
$ws.Cells.Item(122, 13).Value = -250
$ws.Cells.Item(122, 14).Value = -11874.25
# GSM row 132
$ws.Cells.Item(132, 8).Value = 2826.1836
$ws.Cells.Item(132, 9).Value = 2049.3513
$ws.Cells.Item(132, 10).Value = 5221.4165
$ws.Cells.Item(132, 11).Value = 6148.053899999999
$ws.Cells.Item(132, 12).Value = 15664.2495
$ws.Cells.Item(132, 13).Value = -3618.053899999999
$ws.Cells.Item(132, 14).Value = -20724.2495

$ws = $wb.Worksheets.Item("LTW")
# LTW row 14
$ws.Cells.Item(14, 8).Value = 18600
$ws.Cells.Item(14, 10).Value = 2900
$ws.Cells.Item(14, 12).Value = 2900
$ws.Cells.Item(14, 14).Value = -3244
# LTW row 24
$ws.Cells.Item(24, 8).Value = 4000
$ws.Cells.Item(24, 10).Value = 4000
$ws.Cells.Item(24, 12).Value = 4000
$ws.Cells.Item(24, 14).Value = -4686
# LTW row 40
$ws.Cells.Item(40, 8).Value = 3499.2942
$ws.Cells.Item(40, 9).Value = 2994
$ws.Cells.Item(40, 11).Value = 2994
$ws.Cells.Item(40, 13).Value = -2858
# LTW row 41
$ws.Cells.Item(41, 8).Value = 0
$ws.Cells.Item(41, 9).Value = 0
$ws.Cells.Item(41, 11).Value = 0
$ws.Cells.Item(41, 13).ClearContents()
# LTW row 93
$ws.Cells.Item(93, 8).Value = 1490.1111
$ws.Cells.Item(93, 9).Value = 1490.1111
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 11).Value = 1490.1111
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 13).Value = -242.1111000000001
$ws.Cells.Item(93, 14).ClearContents()
# LTW row 100
$ws.Cells.Item(100, 8).Value = 2734.724
$ws.Cells.Item(100, 9).Value = 1850.7
$ws.Cells.Item(100, 10).Value = 3200
$ws.Cells.Item(100, 11).Value = 1850.7
$ws.Cells.Item(100, 12).Value = 3200
$ws.Cells.Item(100, 13).Value = -1309.7
$ws.Cells.Item(100, 14).Value = -4282

$ws = $wb.Worksheets.Item("WVR")
# WVR row 49
$ws.Cells.Item(49, 8).Value = 6000
$ws.Cells.Item(49, 9).Value = 0
$ws.Cells.Item(49, 11).Value = 0
$ws.Cells.Item(49, 13).ClearContents()
# WVR row 81
$ws.Cells.Item(81, 8).Value = 3586.4194
$ws.Cells.Item(81, 9).Value = 2325.2666
$ws.Cells.Item(81, 10).Value = 4768.75
$ws.Cells.Item(81, 11).Value = 4650.5332
$ws.Cells.Item(81, 12).Value = 9537.5
$ws.Cells.Item(81, 13).Value = -3589.5332
$ws.Cells.Item(81, 14).Value = -11659.5
# WVR row 84
$ws.Cells.Item(84, 8).Value = 3586.4194
$ws.Cells.Item(84, 9).Value = 2325.2666
$ws.Cells.Item(84, 10).Value = 4768.75
$ws.Cells.Item(84, 11).Value = 23252.666
$ws.Cells.Item(84, 12).Value = 47687.5
$ws.Cells.Item(84, 13).Value = -17948.666
$ws.Cells.Item(84, 14).Value = -58295.5
# WVR row 92
$ws.Cells.Item(92, 8).Value = 0
$ws.Cells.Item(92, 10).Value = 0
$ws.Cells.Item(92, 12).Value = 0
$ws.Cells.Item(92, 14).ClearContents()
# WVR row 93
$ws.Cells.Item(93, 8).Value = 0
$ws.Cells.Item(93, 10).Value = 0
$ws.Cells.Item(93, 12).Value = 0
$ws.Cells.Item(93, 14).ClearContents()
# WVR row 102
$ws.Cells.Item(102, 8).Value = 0
$ws.Cells.Item(102, 10).Value = 0
$ws.Cells.Item(102, 12).Value = 0
$ws.Cells.Item(102, 14).ClearContents()
# WVR row 136
$ws.Cells.Item(136, 8).Value = 3055.898
$ws.Cells.Item(136, 9).Value = 1220.8529
$ws.Cells.Item(136, 10).Value = 7215.3335
$ws.Cells.Item(136, 11).Value = 3662.5587
$ws.Cells.Item(136, 12).Value = 21646.0005
$ws.Cells.Item(136, 13).Value = -1112.5587
$ws.Cells.Item(136, 14).Value = -26746.0005
